# Add a new "Western Canada"-like distribution row for "Canada" just
# below the header/example row (row 2), pushing the existing example
# rows (previously rows 3-12) down to rows 4-13.
#
# Constants used (since the sandbox doesn't expose the Excel enum names):
#   xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3, shifting rows 3:12 down to 4:13.
$ws.Rows.Item(3).Insert()

# Copy the formatting from row 2 (the existing "Aus bus / Western
# Canada / Saskatchewan" example row) onto the newly inserted row 3,
# so the new row keeps the same cell styles used throughout the sheet.
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's data - same otu/area pattern as row 2, but for
# "Canada" (country) instead of "Western Canada" (geographic_area_name).
$ws.Cells.Item(3, 1).Value = "Aus bus"
$ws.Cells.Item(3, 4).Value = "Canada"
$ws.Cells.Item(3, 5).Value = "Saskatchewan"
$ws.Cells.Item(3, 7).Value = 11

# Match the updated selection left behind in the saved workbook.
$ws.Range("E3").Select()
